$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price report (fecha serial 44610) for "Palta" gets prepended
# above the existing entries for this product block. Insert two blank rows
# at 336/337, which pushes the former rows 336-354 down to 338-356.
$ws.Rows.Item(336).Insert()
$ws.Rows.Item(336).Insert()

# Row 336: Primera calidad
$ws.Range("A336").Value = 4
$ws.Range("B336").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C336").Value = "Los Lagos"
$ws.Range("D336").Value = 44610
$ws.Range("E336").Value = 10
$ws.Range("F336").Value = "Fruta"
$ws.Range("G336").Value = 100106
$ws.Range("H336").Value = "Oleaginosos"
$ws.Range("I336").Value = 100106002
$ws.Range("J336").Value = "Palta"
$ws.Range("K336").Value = "Hass"
$ws.Range("L336").Value = "Primera"
$ws.Range("M336").Value = 300
$ws.Range("N336").Value = 4000
$ws.Range("O336").Value = 4000
$ws.Range("P336").Value = 4000
$ws.Range("Q336").Value = "$/kilo (en caja de 17 kilos)"
$ws.Range("R336").Value = "Provincia de Quillota"
$ws.Range("S336").Value = 4000
$ws.Range("T336").Value = 1

# Row 337: Segunda calidad
$ws.Range("A337").Value = 4
$ws.Range("B337").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C337").Value = "Los Lagos"
$ws.Range("D337").Value = 44610
$ws.Range("E337").Value = 10
$ws.Range("F337").Value = "Fruta"
$ws.Range("G337").Value = 100106
$ws.Range("H337").Value = "Oleaginosos"
$ws.Range("I337").Value = 100106002
$ws.Range("J337").Value = "Palta"
$ws.Range("K337").Value = "Hass"
$ws.Range("L337").Value = "Segunda"
$ws.Range("M337").Value = 200
$ws.Range("N337").Value = 3500
$ws.Range("O337").Value = 3500
$ws.Range("P337").Value = 3500
$ws.Range("Q337").Value = "$/kilo (en caja de 17 kilos)"
$ws.Range("R337").Value = "Provincia de Quillota"
$ws.Range("S337").Value = 3500
$ws.Range("T337").Value = 1
